# Updated cryptos list on Thu Jul 11 03:28:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Price, Volume) updates for rows where only D/E (or just E) change.
$priceVolumeUpdates = @{
    2  = @("57.348.52", "  -1.07%  ")
    3  = @("3.069.53",  "  -0.08%  ")
    4  = @($null,        "  -0.01%  ")
    5  = @("522.48",    "  +1.50%  ")
    6  = @("139.91",    "  -0.68%  ")
    7  = @($null,        "  +0.15%  ")
    8  = @("3.067.02",  "  -0.11%  ")
    9  = @("0.438",     "  +0.67%  ")
    10 = @("7.13",      "  -2.63%  ")
    11 = @("0.107",     "  -0.64%  ")
    12 = @("0.387",     "  +2.29%  ")
    13 = @("3.601.52",  "  +0.04%  ")
    14 = @("0.132",     "  +1.88%  ")
    15 = @("25.28",     "  -5.17%  ")
    16 = @("0.0000162", "  -0.77%  ")
    17 = @("57.403.47", "  -0.95%  ")
    18 = @("3.088.83",  "  +0.69%  ")
    19 = @("6.05",      "  -1.70%  ")
    20 = @("12.61",     "  -1.75%  ")
    21 = @("7.89",      "  -2.39%  ")
    22 = @("337.31",    "  +1.21%  ")
    23 = @($null,        "  +0.23%  ")
    24 = @("0.506",     "  +0.88%  ")
    25 = @("66.95",     "  +2.89%  ")
    26 = @("0.168",     "  -2.44%  ")
    27 = @($null,        "  -0.06%  ")
    28 = @("0.0₃0903", "  -0.25%  ")
    29 = @($null,        "  +0.03%  ")
    30 = @("6.34",      "  -1.32%  ")
    31 = @("7.17",      "  +0.32%  ")
    32 = @("1.86",      "  +3.41%  ")
    33 = @("20.78",     "  -0.01%  ")
    34 = @("1.18",      "  -1.83%  ")
    35 = @("157.51",    "  +1.75%  ")
    36 = @("4.56",      "  +0.35%  ")
    37 = @("6.09",      "  +1.65%  ")
    38 = @("25.87",     "  -5.29%  ")
    39 = @("1.24",      "  -2.30%  ")
    40 = @("0.0659",    "  -2.27%  ")
    41 = @("1.57",      "  +13.32%  ")
    42 = @("3.96",      "  +1.89%  ")
    43 = @("0.678",     "  +3.60%  ")
    44 = @("3.111.60",  "  -0.08%  ")
    45 = @("36.67",     "  +0.29%  ")
    46 = @("0.999",     "  -0.12%  ")
    49 = @("0.984",     "  +4.53%  ")
    50 = @("6.03",      "  +1.54%  ")
    51 = @("20.27",     "  -1.11%  ")
}

foreach ($row in $priceVolumeUpdates.Keys) {
    $vals = $priceVolumeUpdates[$row]
    $price = $vals[0]
    $volume = $vals[1]
    if ($null -ne $price) {
        # Force text format so the dotted-thousands price strings
        # (e.g. "57.348.52") are not reinterpreted as numbers/dates.
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $price
    }
    $ws.Cells.Item($row, 5).Value = $volume
}

# Rows 47 and 48 swap: VeChain <-> Maker fully swapped content.
$ws.Cells.Item(47, 2).Value = "Maker"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.267.03"
$ws.Cells.Item(47, 5).Value = "  -0.54%  "

$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0260"
$ws.Cells.Item(48, 5).Value = "  +2.21%  "
